$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "fj"
$ws.Range("C4").Value = "dhd"
$ws.Range("E5").Value = "hdh"
$ws.Range("G6").Value = "jdj"
$ws.Range("H7").Value = "djd"

[void]$ws.Range("H8").Select()
